$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two confirmed-case totals: the "(revised)" parenthetical
# suffix is dropped now that the as-of timing note lives in its own
# dedicated cell (B1) instead of being appended to every total.
$nbsp = [char]0x00A0
$ws.Range("B7").Value = "29$($nbsp)674"
$ws.Range("B22").Value = "61$($nbsp)940"

# Move the active cell selection down one row (E9 -> E10).
$ws.Range("E10").Select()
